$d = $word.ActiveDocument

# The last paragraph in the body (the "A resource on opioid policy: ...
# 'The Lawhern Files' ..." paragraph) gets wrapped in a "_GoBack" bookmark,
# as Word does automatically to mark the last edited location in the document.
$paragraphs = $d.Paragraphs
$lastPara = $paragraphs.Item($paragraphs.Count)
$r = $lastPara.Range

# Exclude the trailing paragraph mark from the bookmark span.
$bmRange = $d.Range($r.Start, $r.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
